# Realestate Update resale numbers 2024-01-22 22:43
# Append a new data row (row 86) to the CityResaleNum sheet with the
# latest resale-price snapshot for 2024-01-22 22:43:36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 86

# Helper to write a text value into a cell while avoiding Excel's
# automatic type inference (which would otherwise turn values like
# "2024-01-22" into a date serial, or "03" into the number 3).
# We briefly force a text number format so the literal string is
# stored as-is, then clear the formatting again so the cell ends up
# with the sheet's normal (unstyled) appearance, matching the other
# data rows.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $newRow 1 "2024-01-22"
Set-TextValue $newRow 2 "22:43:36"
Set-TextValue $newRow 3 "Monday"
Set-TextValue $newRow 4 "03"

$ws.Cells.Item($newRow, 5).Value  = 138488
$ws.Cells.Item($newRow, 6).Value  = 141017
$ws.Cells.Item($newRow, 7).Value  = 171359
$ws.Cells.Item($newRow, 8).Value  = 148633
$ws.Cells.Item($newRow, 9).Value  = -1
$ws.Cells.Item($newRow, 10).Value = 123217
$ws.Cells.Item($newRow, 11).Value = 223601
$ws.Cells.Item($newRow, 12).Value = 255949
$ws.Cells.Item($newRow, 13).Value = 185402
$ws.Cells.Item($newRow, 14).Value = 110289
$ws.Cells.Item($newRow, 15).Value = 41351
$ws.Cells.Item($newRow, 16).Value = 30892
$ws.Cells.Item($newRow, 17).Value = 73624
$ws.Cells.Item($newRow, 18).Value = -1
$ws.Cells.Item($newRow, 19).Value = 42673
$ws.Cells.Item($newRow, 20).Value = -1
